# Uploaded new data for testing
# Add new execution-number rows to several sheets, matching formatting
# of the existing rows (style s="4" on the text column, copied from the
# row directly above via a CopyOrigin insert so the new cell reuses the
# workbook's existing cell style instead of creating a brand-new one).

$wb = $excel.ActiveWorkbook

function Add-ExecRow($ws, [int]$rowNum, [int]$execValue, [string]$text) {
    $ws.Rows.Item($rowNum).Insert(-4121, 0)
    $ws.Range("A" + $rowNum).Value = $execValue
    $ws.Range("B" + $rowNum).Value = $text
}

# Sheets are touched in the same order the original author entered the
# data, so that newly-appended shared-string indices line up with the
# target workbook (Antennapod, droidweight, GNUCASH-1.0.3, growtracker,
# ATimeTracker, then token last/active).

# --- Antennapod-1.6.2.3 (sheet2) ---
$wsAntennapod = $wb.Worksheets.Item("Antennapod-1.6.2.3")
Add-ExecRow $wsAntennapod 14 13 "antennapod/getevent-1"
Add-ExecRow $wsAntennapod 15 14 "antennapod/getevent-2"
Add-ExecRow $wsAntennapod 16 15 "antennapod/getevent-3"
$wsAntennapod.Range("B16").Select()

# --- droidweight (sheet6) ---
$wsDroidweight = $wb.Worksheets.Item("droidweight")
Add-ExecRow $wsDroidweight 7 6 "getevent-1"
Add-ExecRow $wsDroidweight 8 7 "getevent-2"
$wsDroidweight.Range("G9").Select()

# --- GNUCASH-1.0.3 (sheet4) ---
$wsGnucash103 = $wb.Worksheets.Item("GNUCASH-1.0.3")
Add-ExecRow $wsGnucash103 12 11 "gnucash-1.0.3/getevent-1"
Add-ExecRow $wsGnucash103 13 12 "gnucash-1.0.3/getevent-2"
Add-ExecRow $wsGnucash103 14 13 "gnucash-1.0.3/getevent-3"
$wsGnucash103.Range("G7").Select()

# --- growtracker (sheet7) ---
$wsGrowtracker = $wb.Worksheets.Item("growtracker")
Add-ExecRow $wsGrowtracker 7 6 "growtracker-2.3.1/getevent-1"
Add-ExecRow $wsGrowtracker 8 7 "growtracker-2.3.1/getevent-2"
Add-ExecRow $wsGrowtracker 9 8 "growtracker-2.3.1/getevent-3"
Add-ExecRow $wsGrowtracker 10 9 "growtracker-2.3.1/getevent-4"
Add-ExecRow $wsGrowtracker 11 10 "growtracker-2.3.1/getevent-5"
$wsGrowtracker.Columns.Item(2).ColumnWidth = 23.17
$wsGrowtracker.Range("E11").Select()

# --- ATimeTracker-0.20 (sheet3) ---
$wsATimeTracker = $wb.Worksheets.Item("ATimeTracker-0.20")
Add-ExecRow $wsATimeTracker 14 13 "atimetracker/getevent-1"
$wsATimeTracker.Range("B15").Select()

# --- token (sheet5) ---
$wsToken = $wb.Worksheets.Item("token")
Add-ExecRow $wsToken 7 6 "token-2.10/getevent-1"
Add-ExecRow $wsToken 8 7 "token-2.10/getevent-2"
$wsToken.Range("E8").Select()

# token ends up the active sheet, matching the workbook's saved activeTab.
$wsToken.Activate()
